$p = $ppt.ActivePresentation

# The deck had 15 slides; slides 13 and 14 ("3.2 Spring pom.xml") are
# being removed, leaving the former slide 15 ("End of Chapter") as the
# new slide 13.
$p.Slides.Item(14).Delete()
$p.Slides.Item(13).Delete()
